$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "70.067.99"
$ws.Range("E2").Value = "  +0.12%  "

Set-TextValue $ws.Range("D3") "3.551.98"
$ws.Range("E3").Value = "  +0.37%  "

Set-TextValue $ws.Range("D4") "0.998"
$ws.Range("E4").Value = "  -0.11%  "

Set-TextValue $ws.Range("D5") "603.57"
$ws.Range("E5").Value = "  -2.42%  "

Set-TextValue $ws.Range("D6") "197.27"
$ws.Range("E6").Value = "  +6.35%  "

$ws.Range("E7").Value = "  -0.41%  "

$ws.Range("E9").Value = "  -2.20%  "

Set-TextValue $ws.Range("D10") "0.655"
$ws.Range("E10").Value = "  -0.69%  "

Set-TextValue $ws.Range("D11") "54.19"
$ws.Range("E11").Value = "  +0.84%  "

Set-TextValue $ws.Range("D12") "0.0000305"
$ws.Range("E12").Value = "  -1.15%  "

Set-TextValue $ws.Range("D13") "9.59"
$ws.Range("E13").Value = "  +0.18%  "

Set-TextValue $ws.Range("D14") "4.118.06"
$ws.Range("E14").Value = "  +0.43%  "

Set-TextValue $ws.Range("D15") "605.80"
$ws.Range("E15").Value = "  -4.16%  "

$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
Set-TextValue $ws.Range("D16") "19.23"
$ws.Range("E16").Value = "  +1.25%  "

$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
Set-TextValue $ws.Range("D17") "70.213.34"
$ws.Range("E17").Value = "  +0.23%  "

Set-TextValue $ws.Range("D18") "12.74"
$ws.Range("E18").Value = "  -0.73%  "

Set-TextValue $ws.Range("D19") "3.548.26"
$ws.Range("E19").Value = "  +0.10%  "

$ws.Range("E20").Value = "  +0.50%  "

Set-TextValue $ws.Range("D21") "0.999"
$ws.Range("E21").Value = "  +0.29%  "

Set-TextValue $ws.Range("D22") "18.04"
$ws.Range("E22").Value = "  +2.74%  "

Set-TextValue $ws.Range("D23") "5.30"
$ws.Range("E23").Value = "  +6.81%  "

Set-TextValue $ws.Range("D24") "102.92"
$ws.Range("E24").Value = "  -0.52%  "

Set-TextValue $ws.Range("D25") "4.63"
$ws.Range("E25").Value = "  -2.24%  "

Set-TextValue $ws.Range("D26") "3.14"
$ws.Range("E26").Value = "  +3.78%  "

$ws.Range("E27").Value = "  -0.31%  "

Set-TextValue $ws.Range("D28") "9.68"
$ws.Range("E28").Value = "  -1.81%  "

Set-TextValue $ws.Range("D29") "33.91"
$ws.Range("E29").Value = "  -1.95%  "

Set-TextValue $ws.Range("D30") "4.52"
$ws.Range("E30").Value = "  +24.56%  "

Set-TextValue $ws.Range("D31") "7.16"
$ws.Range("E31").Value = "  +1.16%  "

Set-TextValue $ws.Range("D32") "12.71"
$ws.Range("E32").Value = "  +2.56%  "

$ws.Range("E33").Value = "  +0.46%  "

Set-TextValue $ws.Range("D34") "63.47"
$ws.Range("E34").Value = "  -0.77%  "

Set-TextValue $ws.Range("D35") "0.0₃0840"
$ws.Range("E35").Value = "  +8.05%  "

Set-TextValue $ws.Range("D36") "3.782.70"
$ws.Range("E36").Value = "  +7.24%  "

Set-TextValue $ws.Range("D37") "3.09"
$ws.Range("E37").Value = "  -4.15%  "

$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("E39").Value = "  +3.44%  "

Set-TextValue $ws.Range("D40") "0.397"
$ws.Range("E40").Value = "  -0.85%  "

Set-TextValue $ws.Range("D41") "36.91"
$ws.Range("E41").Value = "  -0.62%  "

Set-TextValue $ws.Range("D42") "487.30"
$ws.Range("E42").Value = "  -8.04%  "

$ws.Range("E43").Value = "  -0.69%  "

$ws.Range("E44").Value = "  -1.20%  "

$ws.Range("E45").Value = "  -3.75%  "

$ws.Range("E46").Value = "  -2.37%  "

Set-TextValue $ws.Range("D47") "3.30"
$ws.Range("E47").Value = "  -1.85%  "

$ws.Range("E48").Value = "  +0.09%  "

Set-TextValue $ws.Range("D49") "8.70"
$ws.Range("E49").Value = "  -3.99%  "

Set-TextValue $ws.Range("D50") "0.000250"
$ws.Range("E50").Value = "  +3.56%  "

Set-TextValue $ws.Range("D51") "131.30"
$ws.Range("E51").Value = "  -1.99%  "
